# Update cryptos list with latest prices and 1h volume changes
# Also fixes row order for WrappedEther / ShibaInu (rows 17-18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.963.72"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.587.26"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.599.24"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.40"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.331"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.042.45"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.914.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.603.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000132"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0716"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.63"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.598"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "267.94"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0953"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0513"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.958.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0220"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.91%  "
